$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.339.17'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '3.810.80'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '700.76'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.06'
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("D7").Value = '3.810.16'
$ws.Range("E7").Value = '  -0.85%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.50'
$ws.Range("E11").Value = '  +2.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.481'
$ws.Range("E12").Value = '  +5.08%  '
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("E14").Value = '  -1.81%  '
$ws.Range("D15").Value = '4.454.33'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").Value = '3.803.09'
$ws.Range("E16").Value = '  -4.52%  '
$ws.Range("D17").Value = '71.394.48'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.22'
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.114'
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '514.82'
$ws.Range("E21").Value = '  +3.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.47'
$ws.Range("E22").Value = '  -1.36%  '
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("E25").Value = '  -3.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.48'
$ws.Range("E26").Value = '  +2.43%  '
$ws.Range("D27").Value = '3.959.02'
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.38'
$ws.Range("E28").Value = '  -2.53%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  -3.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.02'
$ws.Range("E31").Value = '  -5.71%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.23'
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.35'
$ws.Range("E33").Value = '  -2.31%  '
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.19'
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").Value = '3.774.05'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("E39").Value = '  -2.28%  '
$ws.Range("E41").Value = '  +2.48%  '
$ws.Range("E42").Value = '  -1.67%  '
$ws.Range("E43").Value = '  -1.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '172.09'
$ws.Range("E45").Value = '  +5.27%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000309'
$ws.Range("E47").Value = '  -2.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.84'
$ws.Range("E48").Value = '  +2.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '425.53'
$ws.Range("E49").Value = '  +2.45%  '
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("E51").Value = '  -0.51%  '
